# Add "Vapor Facings" KPI column to the Conversion Table sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conversion Table")

# New header in D1, mirroring the existing "Smokeless Facings" column (C).
$ws.Range("D1").Value = "Vapor Facings"

# New values in D2:D5 duplicate the Smokeless Facings values in C2:C5.
$ws.Range("D2").Value = 4
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 16

# Make this sheet the active one and select the cell past the new data,
# matching the workbook's new "current view" state.
$ws.Activate()
$ws.Range("F5").Select()
